$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Направете проучване и имплементирайте алгоритъма Bucket Sort"
$ws.Range("A3").Value = "Направете проучване и имплементирайте алгоритъма Quick Sort"
$ws.Range("A4").Value = "Направете програма, която приема изречение и подрежда по азбучен ред всички думи от него"
$ws.Range("A5").Value = "Направете програма, която приема 10 числа. Подредете във възходящ ред всички четни числа в масив. Подредете всички нечетни числа в нискодящ ред в друг масив"
$ws.Range("A6").Value = "Направете програма, която приема 30 числа и намира всички тройки еднакви числа "

$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 30

$ws.Range("A5").Select()
